# Updated cryptos list on Sat Jun  3 19:57:08 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.153.35"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "'1.891.86"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'306.98"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "'0.5214"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").Value = "'0.3752"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").Value = "'0.07260"
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "'21.15"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "'0.8980"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "'0.08201"
$ws.Range("E12").Value = "  +6.93%  "
$ws.Range("D13").Value = "'96.62"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Value = "'1.892.97"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "'5.270"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "'27.189.51"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").Value = "'5.083"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'10.70"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").Value = "'6.398"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").Value = "'147.58"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "'2.292"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").Value = "'18.18"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "'1.727"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'114.90"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'4.905"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").Value = "'4.790"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").Value = "'0.09227"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "'0.05043"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "'0.7896"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("D34").Value = "'1.216"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").Value = "'3.435"
$ws.Range("E35").Value = "  +4.04%  "
$ws.Range("D36").Value = "'2.974"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "'2.572"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").Value = "'0.5663"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "'1.072"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").Value = "'8.962"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "'6.547"
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").Value = "'115.47"
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("D44").Value = "'0.1517"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "'0.4855"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "'10.09"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "'1.620"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").Value = "'63.25"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").Value = "'0.05937"
$ws.Range("E51").Value = "  -0.04%  "
